$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 892, shifting the existing data
# (old rows 892-985) down to rows 894-987.
$ws.Range("A892:R893").EntireRow.Insert()

# Populate new row 892
$ws.Range("A892").Value = 8
$ws.Range("B892").Value = "Terminal La Palmera de La Serena"
$ws.Range("C892").Value = "Coquimbo"
$ws.Range("D892").Value = 45194
$ws.Range("E892").Value = 4
$ws.Range("F892").Value = 100112043
$ws.Range("G892").Value = "Pepino ensalada"
$ws.Range("H892").Value = "Sin especificar"
$ws.Range("I892").Value = "Primera"
$ws.Range("J892").Value = 600
$ws.Range("K892").Value = 11000
$ws.Range("L892").Value = 12000
$ws.Range("M892").Value = 11500
$ws.Range("N892").Value = "`$/caja 60 unidades"
$ws.Range("O892").Value = "Región de Arica y Parinacota"
$ws.Range("P892").Value = 192
$ws.Range("Q892").Value = 60
$ws.Range("R892").Value = "Hortaliza"

# Populate new row 893
$ws.Range("A893").Value = 8
$ws.Range("B893").Value = "Terminal La Palmera de La Serena"
$ws.Range("C893").Value = "Coquimbo"
$ws.Range("D893").Value = 45194
$ws.Range("E893").Value = 4
$ws.Range("F893").Value = 100112043
$ws.Range("G893").Value = "Pepino ensalada"
$ws.Range("H893").Value = "Sin especificar"
$ws.Range("I893").Value = "Segunda"
$ws.Range("J893").Value = 400
$ws.Range("K893").Value = 7000
$ws.Range("L893").Value = 8000
$ws.Range("M893").Value = 7500
$ws.Range("N893").Value = "`$/caja 80 unidades"
$ws.Range("O893").Value = "Región de Arica y Parinacota"
$ws.Range("P893").Value = 94
$ws.Range("Q893").Value = 80
$ws.Range("R893").Value = "Hortaliza"

# Apply the date number format (matches the style used by the rest of
# column D) to the two new date cells.
$ws.Range("D892:D893").NumberFormat = "YYYY-MM-DD HH:MM:SS"
